$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 395 all held the serial date 45178
# (2023-09-09) and should be bumped by one day to 45179 (2023-09-10).
$ws.Range("C2:C395").Value = 45179
